$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number, report week dates) ---
$ws.Range("A8").Characters(21, 2).Text = "49"
$ws.Range("C9").Characters(27, 10).Text = "12/5/2022"
$ws.Range("C9").Characters(47, 9).Text = "12/11/2022"

# --- Crime-stat table value updates (rows 15-29) ---
$ws.Range("N15").Value = -61.538461538461

$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 13
$ws.Range("G16").Value = 9
$ws.Range("H16").Value = 44.444444444444
$ws.Range("I16").Value = 119
$ws.Range("J16").Value = 84
$ws.Range("K16").Value = 41.666666666666
$ws.Range("L16").Value = 20.202020202020
$ws.Range("M16").Value = -4.8
$ws.Range("N16").Value = -85.765550239234

$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 100
$ws.Range("F17").Value = 11
$ws.Range("G17").Value = 8
$ws.Range("H17").Value = 37.5
$ws.Range("I17").Value = 138
$ws.Range("J17").Value = 105
$ws.Range("K17").Value = 31.428571428571
$ws.Range("L17").Value = 55.056179775280
$ws.Range("M17").Value = 228.571428571429
$ws.Range("N17").Value = -45.454545454545

$ws.Range("C18").Value = 8
$ws.Range("E18").Value = 33.333333333333
$ws.Range("F18").Value = 16
$ws.Range("G18").Value = 17
$ws.Range("H18").Value = -5.882352941176
$ws.Range("I18").Value = 137
$ws.Range("J18").Value = 129
$ws.Range("K18").Value = 6.201550387596
$ws.Range("L18").Value = 24.545454545454
$ws.Range("M18").Value = -18.452380952381
$ws.Range("N18").Value = -81.902245706737

$ws.Range("C19").Value = 18
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = 157.142857142857
$ws.Range("F19").Value = 55
$ws.Range("G19").Value = 38
$ws.Range("H19").Value = 44.736842105263
$ws.Range("I19").Value = 615
$ws.Range("J19").Value = 445
$ws.Range("K19").Value = 38.202247191011
$ws.Range("L19").Value = 69.421487603305
$ws.Range("M19").Value = 60.574412532637
$ws.Range("N19").Value = 29.201680672268

$ws.Range("C20").Value = 2
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 7
$ws.Range("G20").Value = 5
$ws.Range("H20").Value = 40
$ws.Range("I20").Value = 120
$ws.Range("J20").Value = 70
$ws.Range("K20").Value = 71.428571428571
$ws.Range("L20").Value = 62.162162162162
$ws.Range("M20").Value = 34.831460674157
$ws.Range("N20").Value = -87.767584097859

$ws.Range("C21").Value = 32
$ws.Range("D21").Value = 17
$ws.Range("E21").Value = 88.235294117647
$ws.Range("F21").Value = 102
$ws.Range("G21").Value = 78
$ws.Range("H21").Value = 30.769230769230
$ws.Range("I21").Value = 1144
$ws.Range("J21").Value = 845
$ws.Range("K21").Value = 35.384615384615
$ws.Range("L21").Value = 53.145917001338
$ws.Range("M21").Value = 40.713407134071
$ws.Range("N21").Value = -65.697151424287

$ws.Range("C22").Value = 1
$ws.Range("D22").Value = 1
$ws.Range("D22").NumberFormat = "#,##0"
$ws.Range("E22").Value = 0
$ws.Range("E22").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("F22").Value = 4
$ws.Range("H22").Value = 300
$ws.Range("I22").Value = 30
$ws.Range("J22").Value = 22
$ws.Range("K22").Value = 36.363636363636
$ws.Range("L22").Value = 3.448275862068
$ws.Range("M22").Value = -3.225806451612

$ws.Range("G23").Value = 2
$ws.Range("H23").Value = -50

$ws.Range("C24").Value = 35
$ws.Range("D24").Value = 26
$ws.Range("E24").Value = 34.615384615384
$ws.Range("F24").Value = 121
$ws.Range("G24").Value = 102
$ws.Range("H24").Value = 18.627450980392
$ws.Range("I24").Value = 1316
$ws.Range("J24").Value = 946
$ws.Range("K24").Value = 39.112050739957
$ws.Range("L24").Value = 53.558926487748
$ws.Range("M24").Value = 52.845528455284

$ws.Range("C25").Value = 6
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 28
$ws.Range("G25").Value = 27
$ws.Range("H25").Value = 3.703703703703
$ws.Range("I25").Value = 299
$ws.Range("J25").Value = 242
$ws.Range("K25").Value = 23.553719008264
$ws.Range("L25").Value = 69.886363636363
$ws.Range("M25").Value = 68.926553672316

$ws.Range("D26").Value = 1
$ws.Range("D26").NumberFormat = "#,##0"
$ws.Range("E26").Value = -100
$ws.Range("E26").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("G26").Value = 2
$ws.Range("J26").Value = 17
$ws.Range("K26").Value = -5.882352941176
$ws.Range("L26").Value = 14.285714285714

$ws.Range("C27").Value = 2
$ws.Range("D27").Value = 2
$ws.Range("F27").Value = 5
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = 25
$ws.Range("I27").Value = 57
$ws.Range("J27").Value = 45
$ws.Range("K27").Value = 26.666666666666
$ws.Range("L27").Value = 90

$ws.Range("N28").Value = -68.75

$ws.Range("N29").Value = -68.75
